$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GHESS")

# Offshore wind (row 15) shareweight: set every year (columns B:AE) from 1 to 0
$ws.Range("B15:AE15").Value = 0

# Match the saved GHESS view state: select B14:AE15 with B14 as the active
# cell, and scroll the sheet so column D is the first visible column.
$ws.Range("B14:AE15").Select()
$excel.ActiveWindow.ScrollColumn = 4

# Restore "About" as the active/selected sheet (it was active before this
# edit and the workbook view should not change because of it).
$wb.Worksheets.Item("About").Activate()
